$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4260477207521092
$ws.Range("C2").Value = 0.1579759249433508
$ws.Range("E2").Value = 0.1128041530278132
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.4559530195518064
$ws.Range("H2").Value = 0.6229416531491978
$ws.Range("I2").Value = 0.6121505561376495
$ws.Range("K2").Value = 0.2518505838001772
$ws.Range("L2").Value = 0.1930693293216734
$ws.Range("N2").Value = 1.388154561695639
$ws.Range("O2").Value = 2.113888136218279
$ws.Range("B3").Value = 0.3878376787550337
$ws.Range("C3").Value = 0.1582256093673209
$ws.Range("E3").Value = 0.112166094823035
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.4599376453950299
$ws.Range("H3").Value = 0.62785140100398
$ws.Range("I3").Value = 0.6182506463570938
$ws.Range("K3").Value = 0.2201895834674303
$ws.Range("L3").Value = 0.18580819219207
$ws.Range("N3").Value = 1.399237856441111
$ws.Range("O3").Value = 2.132413000128949
$ws.Range("B4").Value = 0.3644201416465762
$ws.Range("C4").Value = 0.1584013677976337
$ws.Range("E4").Value = 0.1118325491951637
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.462668222819886
$ws.Range("H4").Value = 0.631098829088792
$ws.Range("I4").Value = 0.6222844136200756
$ws.Range("K4").Value = 0.2006916284566245
$ws.Range("L4").Value = 0.1814472532157509
$ws.Range("N4").Value = 1.406480081467045
$ws.Range("O4").Value = 2.144869888573666
$ws.Range("B5").Value = 0.3548890312758033
$ws.Range("C5").Value = 0.1584786612100153
$ws.Range("E5").Value = 0.1117113016623357
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.4638523558582719
$ws.Range("H5").Value = 0.6324807845875853
$ws.Range("I5").Value = 0.62400072611063
$ws.Range("K5").Value = 0.1927320515328717
$ws.Range("L5").Value = 0.1796947593348222
$ws.Range("N5").Value = 1.40954134745072
$ws.Range("O5").Value = 2.150218472310627
$ws.Range("B6").Value = 0.3533071284444702
$ws.Range("C6").Value = 0.1584918389350669
$ws.Range("E6").Value = 0.1116920559331227
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.4640532922590168
$ws.Range("H6").Value = 0.6327137983490374
$ws.Range("I6").Value = 0.6242900993030069
$ws.Range("K6").Value = 0.1914095411060686
$ws.Range("L6").Value = 0.1794052499686245
$ws.Range("N6").Value = 1.410056315181986
$ws.Range("O6").Value = 2.151123049359214
$ws.Range("B7").Value = 0.3642915532644224
$ws.Range("C7").Value = 0.1584023872109874
$ws.Range("E7").Value = 0.1118308545447348
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.4626839033724437
$ws.Range("H7").Value = 0.6311172292699183
$ws.Range("I7").Value = 0.6223072667446043
$ws.Range("K7").Value = 0.2005843387493513
$ws.Range("L7").Value = 0.1814235186021733
$ws.Range("N7").Value = 1.406520921243327
$ws.Range("O7").Value = 2.144940918795328
$ws.Range("B8").Value = 0.4128643006413029
$ws.Range("C8").Value = 0.1580573712338555
$ws.Range("E8").Value = 0.1125720859319621
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.4572679656722514
$ws.Range("H8").Value = 0.6245862461141272
$ws.Range("I8").Value = 0.6141940499966232
$ws.Range("K8").Value = 0.2409462391084105
$ws.Range("L8").Value = 0.1905455403637859
$ws.Range("N8").Value = 1.391885437347351
$ws.Range("O8").Value = 2.120050877787946
$ws.Range("B9").Value = 0.5084320778462086
$ws.Range("C9").Value = 0.157557865362417
$ws.Range("E9").Value = 0.1144863378555208
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.4489015444428475
$ws.Range("H9").Value = 0.61362380990996
$ws.Range("I9").Value = 0.6005702244997195
$ws.Range("K9").Value = 0.3196141614215833
$ws.Range("L9").Value = 0.2092024335597671
$ws.Range("N9").Value = 1.366648203216023
$ws.Range("O9").Value = 2.079827679649682
$ws.Range("B10").Value = 0.5788066292675751
$ws.Range("C10").Value = 0.1572974188669392
$ws.Range("E10").Value = 0.116172164177506
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.4441303154703604
$ws.Range("H10").Value = 0.6066908643946221
$ws.Range("I10").Value = 0.5919527294900639
$ws.Range("K10").Value = 0.377094122093979
$ws.Range("L10").Value = 0.2233738831629921
$ws.Range("N10").Value = 1.350210715839609
$ws.Range("O10").Value = 2.055505875254426
$ws.Range("B11").Value = 0.6108503397690583
$ws.Range("C11").Value = 0.157201775364058
$ws.Range("E11").Value = 0.1169994781050683
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.4422587609627868
$ws.Range("H11").Value = 0.6037796041147132
$ws.Range("I11").Value = 0.5883342789417192
$ws.Range("K11").Value = 0.4031696088795229
$ws.Range("L11").Value = 0.2299208075674812
$ws.Range("N11").Value = 1.343188399516727
$ws.Range("O11").Value = 2.045576120813138
$ws.Range("B12").Value = 0.6229880722862617
$ws.Range("C12").Value = 0.1571688168330141
$ws.Range("E12").Value = 0.1173214185891283
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.441593056196723
$ws.Range("H12").Value = 0.6027120124909615
$ws.Range("I12").Value = 0.5870074246889487
$ws.Range("K12").Value = 0.4130327801434248
$ws.Range("L12").Value = 0.2324142670541107
$ws.Range("N12").Value = 1.340594590353476
$ws.Range("O12").Value = 2.041979038554928
$ws.Range("B13").Value = 0.620373853177739
$ws.Range("C13").Value = 0.1571757704527599
$ws.Range("E13").Value = 0.1172516985730034
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.4417345141644162
$ws.Range("H13").Value = 0.6029403886211711
$ws.Range("I13").Value = 0.5872912575376148
$ws.Range("K13").Value = 0.4109090712591978
$ws.Range("L13").Value = 0.2318766229506508
$ws.Range("N13").Value = 1.341150306351402
$ws.Range("O13").Value = 2.042746480962307
$ws.Range("B14").Value = 0.6118488534738447
$ws.Range("C14").Value = 0.1571989986557867
$ws.Range("E14").Value = 0.1170257910863981
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.4422031308274583
$ws.Range("H14").Value = 0.6036910747143054
$ws.Range("I14").Value = 0.5882242486037335
$ws.Range("K14").Value = 0.4039812829952893
$ws.Range("L14").Value = 0.230125660543635
$ws.Range("N14").Value = 1.342973695132883
$ws.Range("O14").Value = 2.04527691802987
$ws.Range("B15").Value = 0.6066274727253642
$ws.Range("C15").Value = 0.1572136503841932
$ws.Range("E15").Value = 0.11688854230378
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.4424957744352795
$ws.Range("H15").Value = 0.6041554275051197
$ws.Range("I15").Value = 0.5888013809535018
$ws.Range("K15").Value = 0.3997363548554915
$ws.Range("L15").Value = 0.2290550013053121
$ws.Range("N15").Value = 1.344099087189143
$ws.Range("O15").Value = 2.046848123691944
$ws.Range("B16").Value = 0.5767130291659157
$ws.Range("C16").Value = 0.1573041268337931
$ws.Range("E16").Value = 0.1161193102223521
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.4442586432543081
$ws.Range("H16").Value = 0.6068859996612801
$ws.Range("I16").Value = 0.5921952738260217
$ws.Range("K16").Value = 0.3753885159603954
$ws.Range("L16").Value = 0.2229480321472153
$ws.Range("N16").Value = 1.350678794525074
$ws.Range("O16").Value = 2.056177632621512
$ws.Range("B17").Value = 0.5583685696319378
$ws.Range("C17").Value = 0.1573654626175696
$ws.Range("E17").Value = 0.1156628628497955
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.4454166840110574
$ws.Range("H17").Value = 0.6086232156267712
$ws.Range("I17").Value = 0.5943545815783793
$ws.Range("K17").Value = 0.3604329021139563
$ws.Range("L17").Value = 0.219227187021005
$ws.Range("N17").Value = 1.354831767530975
$ws.Range("O17").Value = 2.062191506926766
$ws.Range("B18").Value = 0.5478201929443571
$ws.Range("C18").Value = 0.1574028920972168
$ws.Range("E18").Value = 0.1154060157525656
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.4461108906556959
$ws.Range("H18").Value = 0.6096452507687289
$ws.Range("I18").Value = 0.5956249553922213
$ws.Range("K18").Value = 0.3518240484627029
$ws.Range("L18").Value = 0.2170965005356891
$ws.Range("N18").Value = 1.35726329510203
$ws.Range("O18").Value = 2.065757300537115
$ws.Range("B19").Value = 0.5442492124399507
$ws.Range("C19").Value = 0.1574159352047353
$ws.Range("E19").Value = 0.1153200299274317
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.4463507678911185
$ws.Range("H19").Value = 0.6099952174320649
$ws.Range("I19").Value = 0.5960599596220852
$ws.Range("K19").Value = 0.3489080933506727
$ws.Range("L19").Value = 0.2163767130608534
$ws.Range("N19").Value = 1.358093929633007
$ws.Range("O19").Value = 2.066982955768466
$ws.Range("B20").Value = 0.5603210775540788
$ws.Range("C20").Value = 0.1573587108970784
$ws.Range("E20").Value = 0.1157108638711577
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.4452904967918059
$ws.Range("H20").Value = 0.6084359230448584
$ws.Range("I20").Value = 0.5941217806486812
$ws.Range("K20").Value = 0.3620256590457132
$ws.Range("L20").Value = 0.2196223007755549
$ws.Range("N20").Value = 1.354385242400546
$ws.Range("O20").Value = 2.061540269355802
$ws.Range("B21").Value = 0.6143527648893041
$ws.Range("C21").Value = 0.157192087696842
$ws.Range("E21").Value = 0.1170919109793083
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.4420643191712941
$ws.Range("H21").Value = 0.6034696348124271
$ws.Range("I21").Value = 0.5879490294089678
$ws.Range("K21").Value = 0.4060164472733447
$ws.Range("L21").Value = 0.2306395741449734
$ws.Range("N21").Value = 1.342436347496637
$ws.Range("O21").Value = 2.04452924098517
$ws.Range("B22").Value = 0.6496854620086481
$ws.Range("C22").Value = 0.157102178121832
$ws.Range("E22").Value = 0.1180449350744581
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.4402065480148352
$ws.Range("H22").Value = 0.6004269300718263
$ws.Range("I22").Value = 0.5841675951117082
$ws.Range("K22").Value = 0.4347022705865129
$ws.Range("L22").Value = 0.2379231777914725
$ws.Range("N22").Value = 1.335008166214585
$ws.Range("O22").Value = 2.034362214245562
$ws.Range("B23").Value = 0.6308261949151586
$ws.Range("C23").Value = 0.1571484350274801
$ws.Range("E23").Value = 0.117531685156429
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.4411751243653441
$ws.Range("H23").Value = 0.6020323134324883
$ws.Range("I23").Value = 0.5861626879215045
$ws.Range("K23").Value = 0.4193982361356632
$ws.Range("L23").Value = 0.234028215778892
$ws.Range("N23").Value = 1.338937875359157
$ws.Range("O23").Value = 2.039701568417897
$ws.Range("B24").Value = 0.5594383545854384
$ws.Range("C24").Value = 0.1573617565998262
$ws.Range("E24").Value = 0.1156891452554447
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.4453474574679248
$ws.Range("H24").Value = 0.6085205254827954
$ws.Range("I24").Value = 0.5942269397439865
$ws.Range("K24").Value = 0.3613056068306832
$ws.Range("L24").Value = 0.2194436434326406
$ws.Range("N24").Value = 1.354586979575032
$ws.Range("O24").Value = 2.061834356400468
$ws.Range("B25").Value = 0.482547879190065
$ws.Range("C25").Value = 0.1576741843154217
$ws.Range("E25").Value = 0.113919308030372
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.4509234124257944
$ws.Range("H25").Value = 0.6163923019279594
$ws.Range("I25").Value = 0.6040112643502056
$ws.Range("K25").Value = 0.2983864221097861
$ws.Range("L25").Value = 0.2040734088432856
$ws.Range("N25").Value = 1.373105524562781
$ws.Range("O25").Value = 2.089790345815814
